# The workbook's six quarterly-snapshot tabs get "resorted": the tab order
# is reversed so the "总计" (grand total) summary sheet leads, followed by
# the quarters from most to least recent, ending with "2020-Q4".
#
# Target tab order: 总计, 2022-Q2, 2022-Q1, 2021-Q3, 2021-Q2, 2020-Q4

$wb = $excel.ActiveWorkbook

$previouslyActiveSheetName = $wb.ActiveSheet.Name

$wb.Worksheets.Item("总计").Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item("2022-Q2").Move($wb.Worksheets.Item(2))
$wb.Worksheets.Item("2022-Q1").Move($wb.Worksheets.Item(3))
$wb.Worksheets.Item("2021-Q3").Move($wb.Worksheets.Item(4))
$wb.Worksheets.Item("2021-Q2").Move($wb.Worksheets.Item(5))
# "2020-Q4" naturally lands in the last (6th) position once everything
# else has been placed ahead of it.

# Re-select whichever sheet was active before the reorder, so the tab
# that was open for the user stays open (its position just moved).
$wb.Worksheets.Item($previouslyActiveSheetName).Activate()
